$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 18.998697916666668

# --- H10/H11: value change + font color (adds style index) ---
$ws.Range("H10").Value = -3.96
$ws.Range("H10").Font.Color = 0
$ws.Range("H11").Value = -3.96
$ws.Range("H11").Font.Color = 0

# --- H12/H13/H14: value change only (already styled) ---
$ws.Range("H12").Value = -3.96
$ws.Range("H13").Value = -3.96
$ws.Range("H14").Value = -3.96

# --- New row 16 ---
$ws.Range("A16").Value = "2024-09-27_B_e"
$ws.Range("B16").Value = "freq"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0.00001
$ws.Range("E16").Value = 18.2
$ws.Range("F16").Value = 0.537
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -5.1
$ws.Range("I16").Value = -4.9
$ws.Range("J16").Value = 202.14
$ws.Range("K16").Value = 47.2227
$ws.Range("L16").Value = "square"
$ws.Range("M16").Value = 1
$ws.Range("P16").Value = 0
$ws.Range("A16:M16").Font.Color = 0
$ws.Range("P16").Font.Color = 0

# --- New row 17 ---
$ws.Range("A17").Value = "2024-09-27_C_e"
$ws.Range("B17").Value = "freq"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0.00002
$ws.Range("E17").Value = 18.2
$ws.Range("F17").Value = 0.537
$ws.Range("G17").Value = 3.5
$ws.Range("H17").Value = -5.1
$ws.Range("I17").Value = -4.9
$ws.Range("J17").Value = 202.14
$ws.Range("K17").Value = 47.2227
$ws.Range("L17").Value = "square"
$ws.Range("M17").Value = 1
$ws.Range("P17").Value = 0
$ws.Range("A17:M17").Font.Color = 0
$ws.Range("P17").Font.Color = 0

# --- Selection ---
$ws.Range("M17").Select() | Out-Null
